$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 12.5
$ws.Range("B3").Value = 12.5
$ws.Range("C3").Value = 12.3
$ws.Range("C17").Value = 5.6
